# CasosColombia.xlsx — "Add files via upload" update
#
# 1) Twelve cells in existing rows flip between the shared-string "NaN"
#    marker and a real numeric value (both directions), and
# 2) a brand-new data row (177, date 2020-08-28) is appended at the
#    bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix up individual cells that were previously "NaN" placeholders
#        (or, in two cases, the reverse: a stray number becomes "NaN"). ---

# Row 11 / Row 12: BU11, AS12, BU12 had "NaN" -> now real counts of 1.
$ws.Range("BU11").Value = 1
$ws.Range("AS12").Value = 1
$ws.Range("BU12").Value = 1

# Row 18 / Row 19: L18, L19 had a placeholder 1 -> now genuinely "NaN".
$ws.Range("L18").Value = "NaN"
$ws.Range("L19").Value = "NaN"

# Rows 34-37: AK34..AK37 had "NaN" -> now real counts.
$ws.Range("AK34").Value = 1
$ws.Range("AK35").Value = 1
$ws.Range("AK36").Value = 2
$ws.Range("AK37").Value = 2

# Rows 62 / 68: AW62, AW68 had stray numbers -> now "NaN".
$ws.Range("AW62").Value = "NaN"
$ws.Range("AW68").Value = "NaN"

# Row 96: CF96 had "NaN" -> now a real count of 6.
$ws.Range("CF96").Value = 6

# --- 2) Append the new row 177 (date serial 44071 = 2020-08-28). ---

$ws.Range("A177").Value = 44071

$rowVals = @(590520,2696,77289,63758,204065,24804,3367,2707,5536,4740,9655,3697,19113,21437,4779,4282,11942,7232,13709,10900,2736,1031,5612,16763,11365,6485,45667,986,175,251,444,92,49,239,1947,2850,35927,6521,2399,35894,893,20145,1435,7618,1462,1554,3880,1555,927,2464,2586,44781,11891,2351,7340,3649,278,1389,2574,729,1993,8083,8148,8020,13730,1864,813,6746,6021,7144,1439,1442,2611,3023,754,4141,2372,1247,675,1936,1821,1144,896,4559,1289,1133,1178,1482,1381,1509,1103,1039,1081,582,2960,944,789,725,1262,1116,603,715,868,1145,938,1084,839,312,331,662,569,397,530,322,574,699,508,473,357,512,114938,248589,9623,107757,67813,27712,8364)

for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(177, 2 + $i).Value = $rowVals[$i]
}

# Match the saved selection/active-cell state from the edit (bottom-right
# pane's active cell moves to the newly-added last row).
$null = $ws.Range("DX177").Select()
